$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row rework -------------------------------------------------
# Old header (A1:J1): task_folder_name, task_title, project_id, profile_id,
# Task_role, task_mediatype, task_filename, task_filepath, createdDate,
# modifiedDate.
# New header (A1:H1): project_name, assignNameOftaggers,
# numOfItemAssignToTagger, assignNameOfReviewers,
# numOfItemAssignToReviewer, task_mediatype, createdDate, modifiedDate.
#
# Columns I and J go away entirely (used range shrinks to A:H), so clear
# them first -- that also drops "createdDate"/"modifiedDate" from their old
# slots before they get re-introduced at G1/H1.
$ws.Range("I1:J1").ClearContents()

# Write the brand-new header strings. Column F keeps "task_mediatype" as-is
# (same text, same column) so it is left untouched.
$ws.Range("A1").Value = "project_name"
$ws.Range("C1").Value = "numOfItemAssignToTagger"
$ws.Range("E1").Value = "numOfItemAssignToReviewer"
$ws.Range("B1").Value = "assignNameOftaggers"
$ws.Range("D1").Value = "assignNameOfReviewers"

# Re-seat createdDate/modifiedDate at their new columns (G, H).
$ws.Range("G1").Value = "createdDate"
$ws.Range("H1").Value = "modifiedDate"

# --- New sample rows -----------------------------------------------------
$ws.Range("C2").Value = 0
$ws.Range("E2").Value = 0
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0

$ws.Range("C3").Value = 0
$ws.Range("E3").Value = 0
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 0

# --- Column widths to fit the new headers --------------------------------
$ws.Columns.Item(1).ColumnWidth = 12.59
$ws.Columns.Item(2).ColumnWidth = 19.74
$ws.Columns.Item(3).ColumnWidth = 24.74
$ws.Columns.Item(4).ColumnWidth = 22.74
$ws.Columns.Item(5).ColumnWidth = 27.31
$ws.Columns.Item(6).ColumnWidth = 14.45
$ws.Columns.Item(7).ColumnWidth = 11.02
$ws.Columns.Item(8).ColumnWidth = 12.59

# --- Selection moves to the reviewer-count entry area --------------------
$ws.Range("F2:F3").Select() | Out-Null
